# Updated cryptos list values (Price and Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.742.37"
Set-TextValue $ws.Range("E2") "  -2.18%  "

Set-TextValue $ws.Range("D3") "1.796.54"
Set-TextValue $ws.Range("E3") "  -1.67%  "

Set-TextValue $ws.Range("E4") "  -0.10%  "

Set-TextValue $ws.Range("D5") "308.23"
Set-TextValue $ws.Range("E5") "  -1.77%  "

Set-TextValue $ws.Range("D6") "1.001"
Set-TextValue $ws.Range("E6") "  +0.00%  "

Set-TextValue $ws.Range("E7") "  +2.28%  "

Set-TextValue $ws.Range("D8") "0.3711"
Set-TextValue $ws.Range("E8") "  -1.33%  "

Set-TextValue $ws.Range("D9") "0.07250"
Set-TextValue $ws.Range("E9") "  -3.46%  "

Set-TextValue $ws.Range("D10") "0.8542"
Set-TextValue $ws.Range("E10") "  -4.46%  "

Set-TextValue $ws.Range("D11") "20.35"
Set-TextValue $ws.Range("E11") "  -3.04%  "

Set-TextValue $ws.Range("D12") "1.817.28"
Set-TextValue $ws.Range("E12") "  -0.46%  "

Set-TextValue $ws.Range("D13") "5.301"
Set-TextValue $ws.Range("E13") "  -1.81%  "

Set-TextValue $ws.Range("D14") "0.07032"
Set-TextValue $ws.Range("E14") "  -1.13%  "

Set-TextValue $ws.Range("D15") "6.476"
Set-TextValue $ws.Range("E15") "  -4.05%  "

Set-TextValue $ws.Range("D16") "90.36"
Set-TextValue $ws.Range("E16") "  -4.35%  "

Set-TextValue $ws.Range("D17") "1.002"
Set-TextValue $ws.Range("E17") "  -0.11%  "

Set-TextValue $ws.Range("D18") "0.000008626"
Set-TextValue $ws.Range("E18") "  -2.14%  "

Set-TextValue $ws.Range("E19") "  +0.03%  "

Set-TextValue $ws.Range("D20") "14.61"
Set-TextValue $ws.Range("E20") "  -3.95%  "

Set-TextValue $ws.Range("D21") "26.754.30"
Set-TextValue $ws.Range("E21") "  -2.17%  "

Set-TextValue $ws.Range("D22") "5.285"
Set-TextValue $ws.Range("E22") "  +0.12%  "

Set-TextValue $ws.Range("D23") "10.60"
Set-TextValue $ws.Range("E23") "  -3.01%  "

Set-TextValue $ws.Range("D24") "2.036.38"
Set-TextValue $ws.Range("E24") "  -0.89%  "

Set-TextValue $ws.Range("D25") "1.908"
Set-TextValue $ws.Range("E25") "  -4.89%  "

Set-TextValue $ws.Range("D26") "149.47"
Set-TextValue $ws.Range("E26") "  -1.33%  "

Set-TextValue $ws.Range("D29") "5.206"
Set-TextValue $ws.Range("E29") "  -2.76%  "

Set-TextValue $ws.Range("D30") "114.24"
Set-TextValue $ws.Range("E30") "  -3.17%  "

Set-TextValue $ws.Range("D31") "0.08827"
Set-TextValue $ws.Range("E31") "  -0.21%  "

Set-TextValue $ws.Range("D32") "0.7519"
Set-TextValue $ws.Range("E32") "  -3.64%  "

Set-TextValue $ws.Range("E33") "  -3.69%  "

Set-TextValue $ws.Range("D34") "4.427"
Set-TextValue $ws.Range("E34") "  -0.48%  "

Set-TextValue $ws.Range("E35") "  -0.39%  "

Set-TextValue $ws.Range("E36") "  -0.03%  "

Set-TextValue $ws.Range("D37") "1.113"
Set-TextValue $ws.Range("E37") "  +0.60%  "

Set-TextValue $ws.Range("D38") "0.01938"
Set-TextValue $ws.Range("E38") "  -2.54%  "

Set-TextValue $ws.Range("D39") "0.05207"
Set-TextValue $ws.Range("E39") "  -2.20%  "

Set-TextValue $ws.Range("D40") "2.901"
Set-TextValue $ws.Range("E40") "  +0.83%  "

Set-TextValue $ws.Range("D41") "7.146"
Set-TextValue $ws.Range("E41") "  -3.33%  "

Set-TextValue $ws.Range("D42") "2.351"
Set-TextValue $ws.Range("E42") "  +3.54%  "

Set-TextValue $ws.Range("D43") "0.5210"
Set-TextValue $ws.Range("E43") "  -2.14%  "

Set-TextValue $ws.Range("D44") "0.1640"
Set-TextValue $ws.Range("E44") "  -5.03%  "

Set-TextValue $ws.Range("D45") "8.465"
Set-TextValue $ws.Range("E45") "  -3.86%  "

Set-TextValue $ws.Range("D46") "0.4988"
Set-TextValue $ws.Range("E46") "  -3.19%  "

Set-TextValue $ws.Range("D47") "10.24"
Set-TextValue $ws.Range("E47") "  -5.19%  "

Set-TextValue $ws.Range("D48") "104.03"
Set-TextValue $ws.Range("E48") "  -2.11%  "

Set-TextValue $ws.Range("E49") "  -0.04%  "

Set-TextValue $ws.Range("D50") "1.640"
Set-TextValue $ws.Range("E50") "  -3.76%  "

Set-TextValue $ws.Range("D51") "0.06302"
Set-TextValue $ws.Range("E51") "  -1.20%  "

# Rows 27 and 28 swap coin data (LidoDAOToken and EthereumClassic swap positions)
Set-TextValue $ws.Range("B27") "LidoDAOToken"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D27") "2.145"
Set-TextValue $ws.Range("E27") "  -13.79%  "

Set-TextValue $ws.Range("B28") "EthereumClassic"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D28") "18.13"
Set-TextValue $ws.Range("E28") "  -2.69%  "
